# Insert a new weekly price row at row 9 (pushing existing data rows 9-76
# down to 10-77), matching the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9..76 down to 10..77, leaving a blank row 9 to populate.
$ws.Rows.Item(9).Insert()

# Populate the new weekly record in row 9.
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 45022
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112035
$ws.Range("G9").Value = "Bruselas (repollito)"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 15
$ws.Range("K9").Value = 27000
$ws.Range("L9").Value = 27000
$ws.Range("M9").Value = 27000
$ws.Range("N9").Value = "$/malla 17 kilos"
$ws.Range("O9").Value = "Provincia de Quillota"
$ws.Range("P9").Value = 1588
$ws.Range("Q9").Value = 17
$ws.Range("R9").Value = "Hortaliza"
